# Loan RBI, Variable Instalments
# On the "Repayment Schedule" sheet a new (blank) column is inserted
# immediately before column N, pushing the existing "Late" (old N) and
# "Outstanding" (old P) columns one place to the right (to O and Q
# respectively), and the current selection is moved to R8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new blank column before column N - this shifts the existing
# N..P columns (and their formatting) one column to the right.
$ws.Columns("N").Insert()

# Match the new active-cell selection recorded for this sheet.
$ws.Range("R8").Select() | Out-Null
